# Applies the "Add files for different work weeks" update to the
# GS170 Action Tracker workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action Tracker")

# --- Semester goal / strategy narrative (replaces the placeholder prompts) ---
$ws.Range("G3").Value = "My goal by the end of the semester will be to get my first web and computer programming certificate and apply for 5 job offers like a Front-End developer."
$ws.Range("G6").Value = "To achieve this goal, I will strengthen my networking and develop my presentation portfolio.                                                                                This is my me in 30 seconds: My name is Marcos Uc and I am a software engineer and a person passionate about technology, I skillfully develop myself in the Front-End area of all kinds of small or large pages, I really like learning about new technologies and I am constantly improving my skills in web development by taking online courses, I am a proactive person with high standards of work ethics, I like challenges and I am willing to give my best in the work area."

# --- Progress checklist (rows 14-18) moves from "Not Started" to done/in-progress ---
$ws.Range("E14").Value = "Complete"
$ws.Range("E15").Value = "Complete"
$ws.Range("E16").Value = "Complete"
$ws.Range("E17").Value = "In Progress"
$ws.Range("E18").Value = "In Progress"

# --- Week 02 (row 23): progress notes toward semester goal ---
$ws.Range("J23").Value = "Currently I am aplying for job oportunities"
$ws.Range("K23").Value = "Medium"
$ws.Range("L23").Value = "Medium"

# --- Week 05 (row 26): weekly contacts ---
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = "Medium"

# --- Week 06 (row 27): weekly contacts + job opportunities ---
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = "Medium"
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = "High"
$ws.Range("L27").Value = "Medium"

# --- Week 08 (row 29): weekly contacts + job opportunities ---
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = "Low"
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = "High"

# --- Week 09 (row 30): weekly contacts ---
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = "High"

# --- Week 10 (row 31): weekly contacts ---
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = "Medium"

# --- View state: scrolled/selected further down the sheet ---
$ws.Activate()
$ws.Range("H28").Select()
$excel.ActiveWindow.ScrollRow = 13

$wb.Application.Calculate()
